$wb = $excel.ActiveWorkbook

# Rename the existing "data" sheet to "Historical Data"
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "Historical Data"

# Append the new row of historical data (2024-09-01 close price)
$lastRow = 298
$dataSheet.Cells.Item($lastRow, 1).Value = 45536
$dataSheet.Cells.Item($lastRow, 2).Value = 5408.419921875

# Match the date formatting used by the rest of column A
$dataSheet.Cells.Item($lastRow, 1).NumberFormat = $dataSheet.Cells.Item($lastRow - 1, 1).NumberFormat()

# Add the two new (currently empty) sheets after the historical data sheet
$betaSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$betaSheet.Name = "Beta Values"
$betaSheet.PageSetup.LeftMargin = 54
$betaSheet.PageSetup.RightMargin = 54
$betaSheet.PageSetup.TopMargin = 72
$betaSheet.PageSetup.BottomMargin = 72
$betaSheet.PageSetup.HeaderMargin = 36
$betaSheet.PageSetup.FooterMargin = 36

$cyclicalitySheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$cyclicalitySheet.Name = "Cyclicality Labels"
$cyclicalitySheet.PageSetup.LeftMargin = 54
$cyclicalitySheet.PageSetup.RightMargin = 54
$cyclicalitySheet.PageSetup.TopMargin = 72
$cyclicalitySheet.PageSetup.BottomMargin = 72
$cyclicalitySheet.PageSetup.HeaderMargin = 36
$cyclicalitySheet.PageSetup.FooterMargin = 36
